# "bulk upload sample file changes"
#
# The Tags column (N) used a comma as the separator between tag values;
# change it to a semicolon for every data row, e.g.
#   "Rozgar Mela, Finance" -> "Rozgar Mela; Finance"
#
# Also nudge a couple of cosmetic view/layout details that came along with
# the author's re-save: the selected cell and the width of column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content change: comma -> semicolon separator in the Tags column ---
$ws.Range("N2:N5").Value = "Rozgar Mela; Finance"

# --- Cosmetic: widen column N slightly ---
$ws.Columns.Item(14).ColumnWidth = 14.665

# --- Cosmetic: scroll right a bit and leave N9 as the active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$excel.Goto($ws.Range("N9"), $true)
